$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The stray "_GoBack" bookmark that sits in the empty paragraph
#    right after the "О предоставлении информации" heading is gone
#    in the edited document -- remove it if present.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Placeholder renames (braces dropped, new names used) across the
#    whole document body -- every "В связи с ..." paragraph uses the
#    same three placeholders.
# ------------------------------------------------------------------

# First occurrence in the doc is missing the curly braces around
# "reason" (a pre-existing typo), handle it after the braced ones.
$rng = $d.Content
$rng.Find.Execute("{reason}", $true, $true, $false, $false, $false, $true, 1, $false, "requestBase", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("reason", $true, $true, $false, $false, $false, $true, 1, $false, "requestBase", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("{number}", $true, $true, $false, $false, $false, $true, 1, $false, "orderNumber", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("{comment}", $true, $true, $false, $false, $false, $true, 1, $false, "criminal", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Word's "_GoBack" bookmark tracks the site of the most recent
#    edit; in the edited document it now sits inside the SECOND
#    "ст.ст. 103, 173" paragraph, right between "ст.ст. " and
#    "103, 173".
# ------------------------------------------------------------------
$occurrence = 0
$searchRange = $d.Content
while ($searchRange.Find.Execute("ст.ст. 103, 173", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)) {
    $occurrence = $occurrence + 1
    if ($occurrence -eq 2) {
        # "ст.ст. " is 7 characters -- collapse the bookmark right after it,
        # i.e. immediately before "103, 173".
        $markPoint = $searchRange.Duplicate
        $markPoint.Start = $searchRange.Start + 7
        $markPoint.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $markPoint)
        break
    }
    $searchRange.Collapse(0)
    $searchRange.End = $d.Content.End
}
